$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Release date: "25 March 2022" -> "30 March 2022", with "30" and
#    " March " ending up as two separate (but identically-formatted) runs,
#    while the trailing "2022" run is left untouched.
# ---------------------------------------------------------------------------
$dateMatch = $d.Content.Duplicate
$dateMatch.Find.Execute("25 March 2022")
$s = $dateMatch.Start

# Temporarily bump the font size of the "2022" run so that it will not be
# silently re-merged with its neighbour while we edit "25 March ".
$r2022 = $d.Range($s + 9, $s + 13)
$r2022.Font.Size = 99

$r25 = $d.Range($s, $s + 2)
$r25.Text = "30"

# Give " March " a distinct size too, forcing it to split away from "30"
# once we are done, then put both back to their original 14pt (sz 28).
$rMarch = $d.Range($s + 2, $s + 9)
$rMarch.Font.Size = 77
$rMarch2 = $d.Range($s + 2, $s + 9)
$rMarch2.Font.Size = 14

$r2022b = $d.Range($s + 9, $s + 13)
$r2022b.Font.Size = 14

# ---------------------------------------------------------------------------
# 2) Drop the proofing-error split around ".prf" after "...image above".
# ---------------------------------------------------------------------------
$t = " shows the location of these options. Once you have done this, quit EuroScope and when you load it again, it will ask you for a .prf (profile) file."
$d.Content.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2)
